$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 26.99883753119549
$ws.Range("C2").Value = 8.855210140476304
$ws.Range("D2").Value = 4.950014309062996
$ws.Range("E2").Value = 9.366722723881587
$ws.Range("F2").Value = 68.30292685251136
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.29205657185642
$ws.Range("L2").Value = 11.32679254332427
$ws.Range("M2").Value = 21.31983020992996
$ws.Range("B3").Value = 26.93403177437454
$ws.Range("C3").Value = 8.667509640607653
$ws.Range("D3").Value = 4.80565097320004
$ws.Range("E3").Value = 9.353242005042816
$ws.Range("F3").Value = 67.57095694723468
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.28563853136761
$ws.Range("L3").Value = 11.36752994719768
$ws.Range("M3").Value = 21.36956271853759
$ws.Range("B4").Value = 26.90388259944218
$ws.Range("C4").Value = 8.555596636862537
$ws.Range("D4").Value = 4.714771103798433
$ws.Range("E4").Value = 9.344773203179127
$ws.Range("F4").Value = 67.1249201157451
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.2817020043935
$ws.Range("L4").Value = 11.39432731145458
$ws.Range("M4").Value = 21.40543528162797
$ws.Range("B5").Value = 26.89402894073325
$ws.Range("C5").Value = 8.510919017923081
$ws.Range("D5").Value = 4.677211216438421
$ws.Range("E5").Value = 9.341273557804998
$ws.Range("F5").Value = 66.94413012048288
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.28009848419236
$ws.Range("L5").Value = 11.40569693194422
$ws.Range("M5").Value = 21.4213929087886
$ws.Range("B6").Value = 26.892539827185
$ws.Range("C6").Value = 8.503558944484228
$ws.Range("D6").Value = 4.67094383667379
$ws.Range("E6").Value = 9.340689502033406
$ws.Range("F6").Value = 66.91417232428975
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.27983224811926
$ws.Range("L6").Value = 11.40761201876172
$ws.Range("M6").Value = 21.42412346553763
$ws.Range("B7").Value = 26.90373985277548
$ws.Range("C7").Value = 8.554990225492517
$ws.Range("D7").Value = 4.714266634326004
$ws.Range("E7").Value = 9.344726202757657
$ws.Range("F7").Value = 67.1224778204444
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.28168037667662
$ws.Range("L7").Value = 11.39447882517557
$ws.Range("M7").Value = 21.40564507281421
$ws.Range("B8").Value = 26.97449698402238
$ws.Range("C8").Value = 8.789852399693308
$ws.Range("D8").Value = 4.900724548217852
$ws.Range("E8").Value = 9.362114114223273
$ws.Range("F8").Value = 68.04991154310419
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.28984252091186
$ws.Range("L8").Value = 11.34046889989908
$ws.Range("M8").Value = 21.33586884454518
$ws.Range("B9").Value = 27.18926813475657
$ws.Range("C9").Value = 9.27294082656744
$ws.Range("D9").Value = 5.247080012837781
$ws.Range("E9").Value = 9.394708957835221
$ws.Range("F9").Value = 69.89018453783559
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.30590026089204
$ws.Range("L9").Value = 11.24868034554144
$ws.Range("M9").Value = 21.24149060556524
$ws.Range("B10").Value = 27.39254815304901
$ws.Range("C10").Value = 9.636371258224372
$ws.Range("D10").Value = 5.487965326197537
$ws.Range("E10").Value = 9.417765620330291
$ws.Range("F10").Value = 71.24811062517445
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.31775597773338
$ws.Range("L10").Value = 11.18980813921342
$ws.Range("M10").Value = 21.19816546756401
$ws.Range("B11").Value = 27.49465878587651
$ws.Range("C11").Value = 9.802490193877354
$ws.Range("D11").Value = 5.59427137880629
$ws.Range("E11").Value = 9.42806680889268
$ws.Range("F11").Value = 71.86558032317392
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.32316779890947
$ws.Range("L11").Value = 11.16487592193097
$ws.Range("M11").Value = 21.18412688121849
$ws.Range("B12").Value = 27.53468636228564
$ws.Range("C12").Value = 9.865425677904579
$ws.Range("D12").Value = 5.634031928620073
$ws.Range("E12").Value = 9.431941066133207
$ws.Range("F12").Value = 72.09923139851944
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.32522023419793
$ws.Range("L12").Value = 11.15569992257175
$ws.Range("M12").Value = 21.17962741773693
$ws.Range("B13").Value = 27.52600563102484
$ws.Range("C13").Value = 9.851871359145756
$ws.Range("D13").Value = 5.625491171586539
$ws.Range("E13").Value = 9.431107851751547
$ws.Range("F13").Value = 72.04892020073628
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.32477806111901
$ws.Range("L13").Value = 11.15766434874583
$ws.Range("M13").Value = 21.18056012082152
$ws.Range("B14").Value = 27.49792476702428
$ws.Range("C14").Value = 9.807667801631878
$ws.Range("D14").Value = 5.597552559715967
$ws.Range("E14").Value = 9.428386076934611
$ws.Range("F14").Value = 71.88480713514173
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.3233365895254
$ws.Range("L14").Value = 11.16411569351641
$ws.Range("M14").Value = 21.18374033489757
$ws.Range("B15").Value = 27.48090080556677
$ws.Range("C15").Value = 9.780593313405596
$ws.Range("D15").Value = 5.580374175072746
$ws.Range("E15").Value = 9.426715460846314
$ws.Range("F15").Value = 71.78425687328121
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.32245406228908
$ws.Range("L15").Value = 11.16810185952449
$ws.Range("M15").Value = 21.18579469056892
$ws.Range("B16").Value = 27.38606680719507
$ws.Range("C16").Value = 9.625523651884238
$ws.Range("D16").Value = 5.480949946971564
$ws.Range("E16").Value = 9.417088673912767
$ws.Range("F16").Value = 71.20774084045061
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.31740271275288
$ws.Range("L16").Value = 11.19147467410635
$ws.Range("M16").Value = 21.19919713880705
$ws.Range("B17").Value = 27.33034042811834
$ws.Range("C17").Value = 9.530535201172407
$ws.Range("D17").Value = 5.419099659470142
$ws.Range("E17").Value = 9.411135215614971
$ws.Range("F17").Value = 70.85390506350284
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.31430893642233
$ws.Range("L17").Value = 11.20628629167099
$ws.Range("M17").Value = 21.2088724095212
$ws.Range("B18").Value = 27.29919691311223
$ws.Range("C18").Value = 9.475979732572746
$ws.Range("D18").Value = 5.383218151834765
$ws.Range("E18").Value = 9.407693282455156
$ws.Range("F18").Value = 70.65037180860408
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.31253120320314
$ws.Range("L18").Value = 11.21497961817064
$ws.Range("M18").Value = 21.21497102568528
$ws.Range("B19").Value = 27.28880907440458
$ws.Range("C19").Value = 9.457524456576463
$ws.Range("D19").Value = 5.371017374245198
$ws.Range("E19").Value = 9.406524858460392
$ws.Range("F19").Value = 70.58146025939085
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.31192957672932
$ws.Range("L19").Value = 11.21795294534452
$ws.Range("M19").Value = 21.21712752355874
$ws.Range("B20").Value = 27.33617871176682
$ws.Range("C20").Value = 9.54063929237627
$ws.Range("D20").Value = 5.425715672132341
$ws.Range("E20").Value = 9.411770796634006
$ws.Range("F20").Value = 70.8915741466373
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.31463809259729
$ws.Range("L20").Value = 11.20469155791006
$ws.Range("M20").Value = 21.20778721868632
$ws.Range("B21").Value = 27.50613608926698
$ws.Range("C21").Value = 9.820651295967172
$ws.Range("D21").Value = 5.605772430462132
$ws.Range("E21").Value = 9.429186246581798
$ws.Range("F21").Value = 71.93301683302118
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.32375989728374
$ws.Range("L21").Value = 11.16221358198017
$ws.Range("M21").Value = 21.18278405710356
$ws.Range("B22").Value = 27.62512947706033
$ws.Range("C22").Value = 10.00379357003603
$ws.Range("D22").Value = 5.720553678965401
$ws.Range("E22").Value = 9.440413497980039
$ws.Range("F22").Value = 72.61260016147482
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.32973976651244
$ws.Range("L22").Value = 11.13599781575587
$ws.Range("M22").Value = 21.17120340558029
$ws.Range("B23").Value = 27.56090516870066
$ws.Range("C23").Value = 9.906061046901856
$ws.Range("D23").Value = 5.659565205970572
$ws.Range("E23").Value = 9.434435340030213
$ws.Range("F23").Value = 72.25003469091536
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.32654638898513
$ws.Range("L23").Value = 11.14984839454324
$ws.Range("M23").Value = 21.17694831830817
$ws.Range("B24").Value = 27.33353643622205
$ws.Range("C24").Value = 9.536071056064698
$ws.Range("D24").Value = 5.422725576852026
$ws.Range("E24").Value = 9.411483510320135
$ws.Range("F24").Value = 70.87454428145514
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.31448927825693
$ws.Range("L24").Value = 11.20541198278487
$ws.Range("M24").Value = 21.20827616350069
$ws.Range("B25").Value = 27.1231050211084
$ws.Range("C25").Value = 9.140387942764203
$ws.Range("D25").Value = 5.155645772767656
$ws.Range("E25").Value = 9.386049453645583
$ws.Range("F25").Value = 69.39077857013342
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.301547929988
$ws.Range("L25").Value = 11.27200423760324
$ws.Range("M25").Value = 21.26246242929763
